$wb = $excel.ActiveWorkbook

# "red" deck: cards were purchased at dsg, so qty_have now covers qty_need
# (qty_buy, the B-C formula column, drops to 0 for these cards).
$wsRed = $wb.Worksheets.Item("red")
$wsRed.Range("C4").Value = 4
$wsRed.Range("C5").Value = 4
$wsRed.Range("C7").Value = 4
$wsRed.Range("C10").Value = 4

# "blue" is no longer the active tab; update its saved selection.
$wsBlue = $wb.Worksheets.Item("blue")
$wsBlue.Activate() | Out-Null
$wsBlue.Range("C2").Select() | Out-Null

# "red" becomes the active tab, with its own updated selection.
$wsRed.Activate() | Out-Null
$wsRed.Range("C11").Select() | Out-Null
